# SR [2022-09-21]: PS.MergeToMain() creds added
#
# 1. "Team" sheet: the "Report version / 1.0" header row (row 4) is removed
#    - the underlying table (Table2) shrinks from A3:G9 to A3:G8. This also
#      orphans the "1.0"/"Report version" shared strings, which reshuffles
#      the shared-string indices referenced elsewhere (e.g. "Planned
#      Objects" sheet) without changing any of that sheet's actual text.
# 2. "Sheet1" is renamed to "PrivelegedUsers" and populated with a small
#    User / MergeMain credentials table (incl. a mailto hyperlink on the
#    second data row), then becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Team: drop the obsolete "Report version" row -----------------------
$wsTeam = $wb.Worksheets.Item("Team")
$wsTeam.Rows.Item(4).Delete()

# --- 2. Sheet1 -> PrivelegedUsers ------------------------------------------
$wsUsers = $wb.Worksheets.Item("Sheet1")
$wsUsers.Name = "PrivelegedUsers"

$wsUsers.Range("A1").Value = "User"
$wsUsers.Range("B1").Value = "MergeMain"

$wsUsers.Range("A2").Value = "sergiy.razumov@gmail.com"
$wsUsers.Range("B2").Value = $true

$wsUsers.Range("A3").Value = "someemail@com"
$wsUsers.Range("B3").Value = $false

# Hyperlink the second credential (mailto), then style it like a hyperlink.
$wsUsers.Hyperlinks.Add($wsUsers.Range("A3"), "mailto:someemail@com")
$wsUsers.Range("A3").Style = "Hyperlink"

# Turn the populated range into a proper table, named like the sheet.
$loUsers = $wsUsers.ListObjects.Add(1, $wsUsers.Range("A1:B3"), 0, 1)
$loUsers.Name = "PrivelegedUsers"

# Restore Team's selection to match the post-edit view state.
$wsTeam.Activate()
$wsTeam.Range("A4:G4").Select()

# Make PrivelegedUsers the active sheet/tab with A4 selected (must be the
# last sheet activated so it ends up as the workbook's active tab).
$wsUsers.Activate()
$wsUsers.Range("A4").Select()
